# Append 5 new event rows (436-440) to the "Tabelle1" sheet, matching the
# rows that already exist (A=Datum, B=Event, C=Location, D=Stadt, E=Link).
#
# Each new row needs:
#   - a date serial in column A
#   - plain text in B/C/D
#   - a hyperlinked URL in column E (text + a real Hyperlinks collection
#     entry, like every other row in this sheet)
# and all of the touched cells must end up with the same "populated row"
# cell style that the rest of the table uses (not the blank-template style
# the empty rows currently carry).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A template row that already has the correct "populated" formatting for
# each column - used purely as a formatting donor via PasteSpecial.
$formatDonorRow = 435

$newRows = @(
    @(436, 45800, "TAGESRAVER FIGHT NIGHT", "Elektroküche", "Köln", "https://www.instagram.com/reel/DJGt_TpgmBn/?igsh=aHdvZzJyazZtbG8y"),
    @(437, 45815, "PLOYZZ ALL NIGHT LONG", "AJZ", "Lüdenscheid", "https://www.instagram.com/reel/DImM3hwMEzR/?igsh=aWxhZTZmeGRoa2M1"),
    @(438, 45815, "HOPSEN", "Klub Kulb", "Düsseldorf", "https://www.instagram.com/area51.techno?igsh=MWI0amhkbHZsN2RneQ=="),
    @(439, 45791, "#MITTWOCHENENDE", "Odonien", "Köln", "https://www.instagram.com/odonien?igsh=Mm8xbmdxenZrYm84"),
    @(440, 45805, "USB", "Prismatic", "Dortmund", "https://www.instagram.com/reel/DJYarB5syJo/?igsh=NGpiYXM0OTNoeTU1")
)

foreach ($row in $newRows) {
    $r = $row[0]
    $date = $row[1]
    $event = $row[2]
    $location = $row[3]
    $city = $row[4]
    $link = $row[5]

    # Date
    $ws.Cells.Item($r, 1).Value = $date

    # Event / Location / Stadt text
    $ws.Cells.Item($r, 2).Value = $event
    $ws.Cells.Item($r, 3).Value = $location
    $ws.Cells.Item($r, 4).Value = $city

    # Bring B:D formatting in line with the rest of the populated table.
    $ws.Range("B$formatDonorRow`:D$formatDonorRow").Copy()
    $ws.Range("B$r`:D$r").PasteSpecial(-4122)

    # Link column: text + real hyperlink (mirrors every other row).
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 5), $link, "", "", $link)

    # Adding the hyperlink restyles the cell; bring it back in line with
    # the rest of the populated "Link" column.
    $ws.Range("E$formatDonorRow").Copy()
    $ws.Cells.Item($r, 5).PasteSpecial(-4122)
}

$ws.Application.CutCopyMode = $false
